$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet (diff: "repayment_20250912_20250912 (2)" -> "... (3)")
$ws.Name = "repayment_20250912_20250912 (3)"

# Helper: write a text value to a cell while avoiding Excel's automatic
# "looks like a number" coercion, and without leaving a lingering
# NumberFormat override on the cell once done.
function Set-TextCell($addr, $text) {
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $text
    $rng.ClearFormats()
}

# --- Row 5 (Aldi Taufik) ---
$ws.Range("D5").Value = 3
Set-TextCell "E5" "893,721.00"
Set-TextCell "G5" "0.56"
$ws.Range("H5").Value = 1.1879999999999999

# --- Row 6 (Axl Wicaksono) ---
$ws.Range("D6").Value = 3
Set-TextCell "E6" "2,341,266.00"
Set-TextCell "G6" "1.58"
$ws.Range("H6").Value = 527

# --- Row 10 (Erlangga Hutama) ---
$ws.Range("D10").Value = 5
Set-TextCell "E10" "2,032,613.00"
Set-TextCell "G10" "1.62"
$ws.Range("H10").Value = 359
Set-TextCell "K10" "0.80"

# --- Talk_time (column H) updates on other rows ---
$ws.Range("H3").Value = 461
$ws.Range("H4").Value = 1.37
$ws.Range("H7").Value = 865
$ws.Range("H8").Value = 565
$ws.Range("H9").Value = 444
$ws.Range("H11").Value = 944
$ws.Range("H12").Value = 719
$ws.Range("H13").Value = 255
$ws.Range("H14").Value = 219
$ws.Range("H15").Value = 4.2640000000000002
$ws.Range("H16").Value = 550
$ws.Range("H17").Value = 1.2210000000000001
